$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# 1. "is rather straight forward government structural design" -> add semicolon
Replace-Text "is rather straight forward government structural design" "is rather straight forward; government structural design"

# 2. "people would choose clinic." -> add "private"
Replace-Text "people would choose clinic." "people would choose private clinic."

# 3. "old or new hospital has reconstructed to innovative design even better" -> "old hospital ... and even better"
Replace-Text "old or new hospital has reconstructed to innovative design even better" "old hospital has reconstructed to innovative design and even better"

# 4. "commercial buildings, with advance" -> "commercial buildings together with advance"
Replace-Text "commercial buildings, with advance" "commercial buildings together with advance"

# 5. "apart from this is lacking of doctors" -> "apart from lacking of doctors"
Replace-Text "apart from this is lacking of doctors" "apart from lacking of doctors"

# 6. "self-own hospital or even clinic." -> "self-own hospital or even private clinic."
Replace-Text "self-own hospital or even clinic." "self-own hospital or even private clinic."

# 7. "and the how the hospital is designed" -> "and how the hospital is designed"
Replace-Text "and the how the hospital is designed" "and how the hospital is designed"

# 8. "still needs to wait. The hospital" -> "still needs to wait for almost one hour. The hospital"
Replace-Text "still needs to wait. The hospital" "still needs to wait for almost one hour. The hospital"

# 9. "waiting for the doctor to see as." -> "waiting for the doctor to see us."
Replace-Text "waiting for the doctor to see as." "waiting for the doctor to see us."

# 10. "main issue of the problem. Grouping all kinds" -> "main issue of the problem is grouping all kinds"
Replace-Text "main issue of the problem. Grouping all kinds" "main issue of the problem is grouping all kinds"

# 11. "in terms of seriousness. I am not the doctor" -> add "and assign to different areas"
Replace-Text "in terms of seriousness. I am not the doctor" "in terms of seriousness and assign to different areas. I am not the doctor"

# 12. "usually with 5 to 8 minutes" -> "usually within 5 to 8 minutes"
Replace-Text "usually with 5 to 8 minutes" "usually within 5 to 8 minutes"

# 13. "not more than 15 minutes. There is another" -> add "to diagnosis a case"
Replace-Text "not more than 15 minutes. There is another" "not more than 15 minutes to diagnosis a case. There is another"

# 14. "As I said everything as a value." -> "As I said everything has a value."
Replace-Text "As I said everything as a value." "As I said everything has a value."

Write-Output "All replacements applied."
